# "fixes to green graphs"
# The four school-year labels that feed the (green) charts get re-typed with a
# two-space indent and 4-digit years ("2018-19" -> "  2018-2019", etc.) across
# all four data blocks on Sheet1. Selection is left on A4, matching the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = @("  2018-2019", "  2019-2020", "  2020-2021", "  2021-2022")

# Column A rows for the four repeated "school year" blocks (EL, IEP, Econ Dis, Homeless)
$blockStartRows = @(4, 12, 20, 28)

foreach ($startRow in $blockStartRows) {
    for ($offset = 0; $offset -lt $years.Length; $offset++) {
        $row = $startRow + $offset
        $ws.Cells.Item($row, 1).Value = $years[$offset]
    }
}

[void]$ws.Range("A4").Select()
